$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Range("D18").Value = "TIMESTAMP"
$ws.Range("D20").Value = "TIMESTAMP"
